# dyes_DFT_PBE_eth_dataset.xlsx - Sheet1 edit
# Drops the Total_Energy_Hartree, Solvation_Energy_eV, Surface_Area_A2 and
# Molecular_Volume_A3 columns (old E:H), which shifts the old
# Max_Absorption_nm / Max_f_osc columns (old I:J) left into E:F, then adds a
# new computed Max_Excitation_eV column in G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old columns E (Total_Energy_Hartree), F (Solvation_Energy_eV),
# G (Surface_Area_A2), H (Molecular_Volume_A3). Old I/J (Max_Absorption_nm /
# Max_f_osc) shift left to become the new E/F, headers included.
$ws.Range("E1:H1").EntireColumn.Delete()

# New third metric column - copy the header formatting (bold, centered,
# bordered) from the neighbouring F1 header cell, then set its text.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Max_Excitation_eV"

$excitation = @{
  2 = 2.13; 3 = 2.22; 4 = 2.43; 5 = 2.09; 6 = 2.35; 7 = 2.34; 8 = 2.21;
  9 = 1.95; 10 = 1.82; 11 = 3.14; 14 = 3.12; 16 = 3.04; 17 = 3.03;
  18 = 3.45; 19 = 2.62; 20 = 2.39; 21 = 2.49; 22 = 2.34; 23 = 2.6;
  24 = 3.21; 25 = 2.35; 26 = 3.25; 27 = 1.55; 28 = 2.46; 29 = 2.36;
  30 = 1.85; 31 = 1.87; 32 = 3.62; 33 = 4.02; 34 = 3.74; 35 = 2.36;
  36 = 2.96; 37 = 2.66; 38 = 2.37; 39 = 2.29; 40 = 2.22; 41 = 2.13
}

foreach ($row in $excitation.Keys) {
  $ws.Cells.Item($row, 7).Value = $excitation[$row]
}

# Rows 12, 13 and 15 have no absorption data, so the new column is left
# untouched there too (stays blank, matching the blank Max_Absorption_nm /
# Max_f_osc already shifted into E/F for those rows).
